$p = $ppt.ActivePresentation

# Slide 2: TextBox "The" + " " + "Moon" -> single run "The Moon"
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(2)
$tb2.TextFrame.TextRange.Text = "temp"
$tb2.TextFrame.TextRange.Text = "The Moon"

# Slide 3: Title "One" + " " + "More" -> single run "One More"
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "temp"
$title3.TextFrame.TextRange.Text = "One More"

# Slide 3: TextBox "The" + " " + "Moon" -> single run "The Moon"
$tb3 = $s3.Shapes.Item(3)
$tb3.TextFrame.TextRange.Text = "temp"
$tb3.TextFrame.TextRange.Text = "The Moon"
